# Add two new report-export columns ("code" and "nextaction") between the
# existing "payment" (I) and "nextkin" (J) columns. Inserting whole columns
# shifts the old J:K ("nextkin"/"kinphone") data over to L:M automatically,
# carrying their values/styles with them - exactly like a user doing
# Insert > Sheet Columns in the Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at J:K (pushes old J:K -> L:M).
$ws.Range("J1:K1").EntireColumn.Insert()

# Header row: label the two new columns. They inherit the "payment" header's
# style (bold font, grey fill, border) from the column to their left.
$ws.Range("J1").Value = "code"
$ws.Range("K1").Value = "nextaction"

# The inserted columns also inherited the thin box border that ran down
# column I; the new "code"/"nextaction" columns are borderless like the
# rest of the blank interior cells, so strip it back off.
$ws.Range("J2:K11").Borders.LineStyle = 0

# Put the active selection on the new J1 cell (matches the saved view).
$null = $ws.Range("J1").Select()
